$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 (shifts current rows 13-23 down to 14-24).
# This new row carries the teacher-name value that belongs next to the
# "Docentes responsaveis:" label now sitting alone in row 12.
$ws.Rows.Item(13).Insert()

# The newly inserted row copies formatting from the row above, which puts
# a stray formatted (but empty) cell in column A - remove it entirely.
$ws.Cells.Item(13, 1).Clear()

# Row 13 (new): value for "Docentes responsaveis:" (row 12's label)
$ws.Cells.Item(13, 2).Value = "5111420 - Talita Martins Lacerda"
$ws.Cells.Item(13, 3).Value = "5111420 - Talita Martins Lacerda"
$ws.Cells.Item(13, 2).Font.Bold = $false
$ws.Cells.Item(13, 2).WrapText = $true
$ws.Cells.Item(13, 2).VerticalAlignment = -4160
$ws.Cells.Item(13, 3).Font.Bold = $false
$ws.Cells.Item(13, 3).WrapText = $true
$ws.Cells.Item(13, 3).VerticalAlignment = -4160
$ws.Cells.Item(13, 3).Font.Color = 255

# Row 10: "Objetivos:" / "Objectives:" body text
$objetivos = "Apresentar os conceitos básicos da ciência dos polímeros, incluindo as reações químicas e os principais métodos de caracterização envolvidos na preparação destes materiais. Estes fundamentos serão usados para introduzir os alunos aos polímeros de fontes renováveis, produzidos a partir de unidades monoméricas extraídas da biomassa."
$ws.Cells.Item(10, 2).Value = $objetivos
$ws.Cells.Item(10, 3).Value = $objetivos

# Row 14: "Programa resumido:" / "Short syllabus:" body text (was "Semestral")
$programaResumido = "Fundamentos sobre a química dos polímeros; Mecanismos de polimerização; Caracterização e propriedades gerais dos polímeros; Monômeros derivados da biomassa e principais polímeros obtidos a partir deles."
$ws.Cells.Item(14, 2).Value = $programaResumido
$ws.Cells.Item(14, 3).Value = $programaResumido

# Row 16: "Programa:" / "Syllabus:" body text
$programa = "Fundamentos sobre a química dos polímeros: composição e estrutura, nomenclatura, polímeros lineares, ramificados, e entrecruzados, massa molar média, propriedades físicas (comportamentos cristalino e amorfo); Mecanismos de polimerização: poliadição e policondensação; Caracterização e propriedades gerais dos polímeros: espectroscopias de infravermelho, FTIR, e ressonância magnética nuclear, RMN, propriedades mecânicas e térmicas; Introdução aos materiais derivados de fontes renováveis; Rotas não-fósseis para a obtenção de monômeros usuais (etileno, propileno, glicerol e derivados); Monômeros exclusivamente obtidos de fontes renováveis (terpenos e terpenóides, monômeros derivados do breu, monômeros derivados de açúcares, ácidos carboxílicos e aminoácidos, furanos, óleos vegetais e derivados); Estado da arte e projeções futuras para os polímeros derivados de fontes renováveis."
$ws.Cells.Item(16, 2).Value = $programa
$ws.Cells.Item(16, 3).Value = $programa

# Row 22: "Bibliografia:" body text
$bibliografia = "Eloisa B. Mano; Introdução a Polímeros, Editora Edgard BlücherLtda, São Paulo, 1999.Elizabete F. Lucas, Bluma G. Soares, Elisabeth E. C. Monteiro; Caracterização de polímeros: determinaçãoo de peso molecular e análise térmica. E-papers Serviços Editoriais Ltda, Rio de Janeiro, 2001.Fred J. Davis; PolymerChemistry: a practical approach. Oxford University Press Inc., New York, 2004.George Odian; Principles of Polymerization.John Wiley and Sons, New Jersey, 2004.Mohamed N. Belgacem, Alessandro Gandini; Monomers, polymers and composites from renewable resources.ElsevierLtda, Amsterdam, 2008."
$ws.Cells.Item(22, 2).Value = $bibliografia
$ws.Cells.Item(22, 3).Value = $bibliografia

# Column A is no longer sharing a <col> definition with column B - split it
# out into its own entry (same visual width as before).
$ws.Columns.Item(1).ColumnWidth = 29.83
